# Issue #57: Make genre required with PBCore controlled vocabulary.
# This fixture (badColumnName_nonRequired.xlsx) gains a new "Genre" column
# (U) in the manifest header row, with a sample value ("Aviation") for the
# first data row, reflecting the new genre-override column supported by
# the batch-ingest manifest spreadsheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column U: header "Genre" (row 2 is the header row in this sheet)
# and a sample controlled-vocabulary value "Aviation" for the first data
# row (row 3). Row 4 is left blank in this column, matching the source
# fixture.
$ws.Cells.Item(2, 21).Value = "Genre"
$ws.Cells.Item(3, 21).Value = "Aviation"

# Match the author's final selection/active cell on the newly added column.
$ws.Range("U3").Select()
